# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1. Remove the two "audit_usefulness" related columns.
#    Column G ("audit_usefulness") is dropped entirely, and the column
#    that was "w_audit_usefulness" (which lands on L after the first
#    delete) is dropped too. Everything to the right shifts left,
#    turning the old A1:O3 range into A1:M3 - matching the new header
#    layout:
#      G -> w_evidence_extraction_quality
#      H -> w_coverage_of_debiasing_mitigation_dimensions
#      I -> w_structure_and_formatting
#      J -> w_relevance_and_faithfulness
#      K -> w_identification_of_missing_disclosures
#      L -> weighted_final_score
#      M -> justification
# -----------------------------------------------------------------------
$ws.Columns("G").Delete()
$ws.Columns("L").Delete()

# -----------------------------------------------------------------------
# 2. Update the scored values / weights / justification text for the two
#    data rows to their new (re-reviewed) figures.
# -----------------------------------------------------------------------

# Row 2 - gpt-4o
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 0.75
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0.4
$ws.Range("J2").Value = 0.3
$ws.Range("K2").Value = 0.9
$ws.Range("L2").Value = 3.35
$ws.Range("M2").Value = "The report provides a structured overview of debiasing strategies with clear sections and tables, earning a 4 in structure and formatting. However, the evidence extraction quality is rated 3 due to a lack of full sentence quotations and some missing validation details. Coverage of debiasing dimensions is strong, with a variety of methods discussed, but some common strategies like resampling are not explicitly documented, leading to a score of 4. Relevance and faithfulness are reasonable, but some claims lack direct support from the sources, resulting in a 3. Missing disclosures are identified, but not comprehensively, also scoring a 3. Overall, the report is useful for audits but could benefit from more detailed validation evidence and explicit documentation of all methods."

# Row 3 - ollama_mistral
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 0.75
$ws.Range("H3").Value = 0.75
$ws.Range("I3").Value = 0.4
$ws.Range("J3").Value = 0.3
$ws.Range("K3").Value = 0.6
$ws.Range("L3").Value = 2.8
$ws.Range("M3").Value = "The report provides a reasonable extraction of evidence with correct citations, but lacks full sentences in some excerpts, affecting evidence extraction quality. Coverage of debiasing methods is broad, yet lacks depth in validation details, particularly for pre-processing and post-processing methods. The structure is clear and well-organized, aiding readability. However, relevance is compromised by some unsupported assumptions, such as the validation status of certain methods. Missing disclosures are not thoroughly identified, especially regarding the absence of validation metrics. The audit usefulness is moderate, as the report is traceable but lacks detailed validation evidence, limiting its utility for comprehensive audits."
